$d = $word.ActiveDocument

# --- Simple text replacements (unique substrings) ---

# Invoice number in the title line
$d.Content.Find.Execute("1/INVOICE", $true, $false, $false, $false, $false, $true, 1, $false, "2021-002/INVOICE", 2) | Out-Null

# "NO. 1  ..." -> "NO. 2021-002  ..." (preserve trailing spaces)
$d.Content.Find.Execute("NO. 1                           ", $true, $false, $false, $false, $false, $true, 1, $false, "NO. 2021-002                           ", 2) | Out-Null

# Seller block
$d.Content.Find.Execute('ООО "Hookah Retrofit"', $true, $false, $false, $false, $false, $true, 1, $false, "ИП Геркен П.В.", 2) | Out-Null
$d.Content.Find.Execute("Address: Some address,", $true, $false, $false, $false, $false, $true, 1, $false, "Address: г. Новочеркасск, Красный спуск, д.6,", 2) | Out-Null

# Buyer block
$d.Content.Find.Execute("Company name: Canada Black Smoke", $true, $false, $false, $false, $false, $true, 1, $false, "Company name: Novochikago Alfa Smoke", 2) | Out-Null
$d.Content.Find.Execute("Company address: Canada", $true, $false, $false, $false, $false, $true, 1, $false, "Company address: Россия", 2) | Out-Null

# Contract number / date
$d.Content.Find.Execute("№/No. 2021-002  от/from 2021-03-23", $true, $false, $false, $false, $false, $true, 1, $false, "№/No. 2021-001  от/from 2021-04-06", 2) | Out-Null

# Loading place
$d.Content.Find.Execute("Loading Place: г. Санкт-Петербург", $true, $false, $false, $false, $false, $true, 1, $false, "Loading Place: Новочеркасск", 2) | Out-Null

# Currency name, all 4 occurrences in the table header row
while ($d.Content.Find.Execute("Датская крона", $true, $false, $false, $false, $false, $true, 1, $false, "Доллар США", 2)) {}

# Shipping marks text, both occurrences
while ($d.Content.Find.Execute("Осторожно!", $true, $false, $false, $false, $false, $true, 1, $false, "Хрупкий груз", 2)) {}

# --- Remove the three goods line-item rows from the table ---
$t = $d.Tables.Item(1)
$t.Rows.Item(5).Delete()
$t.Rows.Item(5).Delete()
$t.Rows.Item(5).Delete()

# --- Update the grand total ---
$d.Content.Find.Execute("60500.00", $true, $false, $false, $false, $false, $true, 1, $false, "0", 2) | Out-Null
